$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6577333333333334
$ws.Range("H2").Value = 1.9732
$ws.Range("M2").Value = 14.65767833333333
$ws.Range("N2").Value = 43.973035
$ws.Range("O2").Value = 0.2345581433878666
$ws.Range("P2").Value = 0.2345581433878665
$ws.Range("Q2").Value = 9.640843629111114
$ws.Range("R2").Value = 86.76759266200001
$ws.Range("S2").Value = 0.2345581433878666
$ws.Range("T2").Value = 0.2345581433878665

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6577333333333334
$ws.Range("H3").Value = 1.9732
$ws.Range("M3").Value = 31.695371
$ws.Range("N3").Value = 95.086113
$ws.Range("O3").Value = 0.5072022462686253
$ws.Range("P3").Value = 0.5072022462686253
$ws.Range("Q3").Value = 20.84710201906667
$ws.Range("R3").Value = 187.6239181716
$ws.Range("S3").Value = 0.5072022462686253
$ws.Range("T3").Value = 0.5072022462686253

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6577333333333334
$ws.Range("H4").Value = 1.9732
$ws.Range("M4").Value = 16.13754733333333
$ws.Range("N4").Value = 48.41264200000001
$ws.Range("O4").Value = 0.2582396103435082
$ws.Range("P4").Value = 0.2582396103435082
$ws.Range("Q4").Value = 10.61420279937778
$ws.Range("R4").Value = 95.52782519440001
$ws.Range("S4").Value = 0.2582396103435082
$ws.Range("T4").Value = 0.2582396103435082
